# Update "Count result" sheet (Sheet1) with the new binning data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (bin index 0): only Count (F2) changes
$ws.Range("F2").Value = 20

# Row 3 (bin index 1)
$ws.Range("B3").Value = "06:30:00"
$ws.Range("C3").Value = "09:00:00"
$ws.Range("D3").Value = 39
$ws.Range("E3").Value = 54
$ws.Range("F3").Value = 400

# Row 4 (bin index 2)
$ws.Range("B4").Value = "09:00:00"
$ws.Range("C4").Value = "12:00:00"
$ws.Range("D4").Value = 54
$ws.Range("E4").Value = 72
$ws.Range("F4").Value = 100

# Row 5 (bin index 3)
$ws.Range("B5").Value = "12:00:00"
$ws.Range("C5").Value = "16:00:00"
$ws.Range("D5").Value = 72
$ws.Range("E5").Value = 96
$ws.Range("F5").Value = 120

# Row 6 (bin index 4)
$ws.Range("B6").Value = "16:00:00"
$ws.Range("C6").Value = "18:00:00"
$ws.Range("D6").Value = 96
$ws.Range("E6").Value = 108
$ws.Range("F6").Value = 200

# Row 7 (bin index 5) -- new row, copy formatting from row 6's A cell (style index 1)
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "18:00:00"
$ws.Range("C7").Value = "22:00:00"
$ws.Range("D7").Value = 108
$ws.Range("E7").Value = 132
$ws.Range("F7").Value = 250

# Row 8 (bin index 6) -- new row, copy formatting from row 6's A cell (style index 1)
$ws.Range("A6").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "22:00:00"
$ws.Range("C8").Value = "00:00:00"
$ws.Range("D8").Value = 132
$ws.Range("E8").Value = 144
$ws.Range("F8").Value = 100

$excel.CutCopyMode = $false
